$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.364465666666667
$ws.Range("H2").Value = 4.093397
$ws.Range("I2").Value = 0.004273090055562101
$ws.Range("J2").Value = 0.004273090055562101
$ws.Range("M2").Value = 16.27546433333333
$ws.Range("N2").Value = 48.826393
$ws.Range("O2").Value = 0.06628560529319844
$ws.Range("P2").Value = 0.06628560529319844
$ws.Range("Q2").Value = 22.20731229189122
$ws.Range("R2").Value = 199.8658106270209
$ws.Range("S2").Value = 0.0002832443608052808
$ws.Range("T2").Value = 0.0002832443608052808
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.364465666666667
$ws.Range("H3").Value = 4.093397
$ws.Range("I3").Value = 0.004273090055562101
$ws.Range("J3").Value = 0.004273090055562101
$ws.Range("O3").Value = 0.3480686258826592
$ws.Range("P3").Value = 0.3480686258826592
$ws.Range("Q3").Value = 116.6115726000439
$ws.Range("R3").Value = 1049.504153400395
$ws.Range("S3").Value = 0.001487328583912356
$ws.Range("T3").Value = 0.001487328583912356
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.364465666666667
$ws.Range("H4").Value = 4.093397
$ws.Range("I4").Value = 0.004273090055562101
$ws.Range("J4").Value = 0.004273090055562101
$ws.Range("M4").Value = 42.61351133333333
$ws.Range("N4").Value = 127.840534
$ws.Range("O4").Value = 0.17355341356458
$ws.Range("P4").Value = 0.17355341356458
$ws.Range("Q4").Value = 58.14467315044421
$ws.Range("R4").Value = 523.3020583539979
$ws.Range("S4").Value = 0.0007416093656116634
$ws.Range("T4").Value = 0.0007416093656116635
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.364465666666667
$ws.Range("H5").Value = 4.093397
$ws.Range("I5").Value = 0.004273090055562101
$ws.Range("J5").Value = 0.004273090055562101
$ws.Range("M5").Value = 101.183272
$ws.Range("N5").Value = 303.549816
$ws.Range("O5").Value = 0.4120923552595624
$ws.Range("P5").Value = 0.4120923552595624
$ws.Range("Q5").Value = 138.0611006849946
$ws.Range("R5").Value = 1242.549906164952
$ws.Range("S5").Value = 0.001760907745232801
$ws.Range("T5").Value = 0.001760907745232801
$ws.Range("I6").Value = 0.9864275102545285
$ws.Range("J6").Value = 0.9864275102545286
$ws.Range("M6").Value = 16.27546433333333
$ws.Range("N6").Value = 48.826393
$ws.Range("O6").Value = 0.06628560529319844
$ws.Range("P6").Value = 0.06628560529319844
$ws.Range("Q6").Value = 5126.478377169012
$ws.Range("R6").Value = 46138.30539452111
$ws.Range("S6").Value = 0.06538594459508412
$ws.Range("T6").Value = 0.06538594459508414
$ws.Range("I7").Value = 0.9864275102545285
$ws.Range("J7").Value = 0.9864275102545286
$ws.Range("O7").Value = 0.3480686258826592
$ws.Range("P7").Value = 0.3480686258826592
$ws.Range("S7").Value = 0.3433444680271464
$ws.Range("T7").Value = 0.3433444680271465
$ws.Range("I8").Value = 0.9864275102545285
$ws.Range("J8").Value = 0.9864275102545286
$ws.Range("M8").Value = 42.61351133333333
$ws.Range("N8").Value = 127.840534
$ws.Range("O8").Value = 0.17355341356458
$ws.Range("P8").Value = 0.17355341356458
$ws.Range("Q8").Value = 13422.48921145455
$ws.Range("R8").Value = 120802.4029030909
$ws.Range("S8").Value = 0.1711978616386831
$ws.Range("T8").Value = 0.1711978616386832
$ws.Range("I9").Value = 0.9864275102545285
$ws.Range("J9").Value = 0.9864275102545286
$ws.Range("M9").Value = 101.183272
$ws.Range("N9").Value = 303.549816
$ws.Range("O9").Value = 0.4120923552595624
$ws.Range("P9").Value = 0.4120923552595624
$ws.Range("Q9").Value = 31870.90981956484
$ws.Range("R9").Value = 286838.1883760836
$ws.Range("S9").Value = 0.4064992359936148
$ws.Range("T9").Value = 0.4064992359936148
$ws.Range("G10").Value = 1.561621333333333
$ws.Range("H10").Value = 4.684864
$ws.Range("I10").Value = 0.004890521434901351
$ws.Range("J10").Value = 0.004890521434901352
$ws.Range("M10").Value = 16.27546433333333
$ws.Range("N10").Value = 48.826393
$ws.Range("O10").Value = 0.06628560529319844
$ws.Range("P10").Value = 0.06628560529319844
$ws.Range("Q10").Value = 25.41611231283911
$ws.Range("R10").Value = 228.745010815552
$ws.Range("S10").Value = 0.0003241711735117974
$ws.Range("T10").Value = 0.0003241711735117975
$ws.Range("G11").Value = 1.561621333333333
$ws.Range("H11").Value = 4.684864
$ws.Range("I11").Value = 0.004890521434901351
$ws.Range("J11").Value = 0.004890521434901352
$ws.Range("O11").Value = 0.3480686258826592
$ws.Range("P11").Value = 0.3480686258826592
$ws.Range("Q11").Value = 133.4611224998045
$ws.Range("R11").Value = 1201.15010249824
$ws.Range("S11").Value = 0.001702237075695804
$ws.Range("T11").Value = 0.001702237075695804
$ws.Range("G12").Value = 1.561621333333333
$ws.Range("H12").Value = 4.684864
$ws.Range("I12").Value = 0.004890521434901351
$ws.Range("J12").Value = 0.004890521434901352
$ws.Range("M12").Value = 42.61351133333333
$ws.Range("N12").Value = 127.840534
$ws.Range("O12").Value = 0.17355341356458
$ws.Range("P12").Value = 0.17355341356458
$ws.Range("Q12").Value = 66.5461683863751
$ws.Range("R12").Value = 598.9155154773759
$ws.Range("S12").Value = 0.0008487666891378773
$ws.Range("T12").Value = 0.0008487666891378777
$ws.Range("G13").Value = 1.561621333333333
$ws.Range("H13").Value = 4.684864
$ws.Range("I13").Value = 0.004890521434901351
$ws.Range("J13").Value = 0.004890521434901352
$ws.Range("M13").Value = 101.183272
$ws.Range("N13").Value = 303.549816
$ws.Range("O13").Value = 0.4120923552595624
$ws.Range("P13").Value = 0.4120923552595624
$ws.Range("Q13").Value = 158.0099561316693
$ws.Range("R13").Value = 1422.089605185024
$ws.Range("S13").Value = 0.002015346496555872
$ws.Range("T13").Value = 0.002015346496555873
$ws.Range("G14").Value = 1.407825
$ws.Range("H14").Value = 4.223475
$ws.Range("I14").Value = 0.00440887825500804
$ws.Range("J14").Value = 0.00440887825500804
$ws.Range("M14").Value = 16.27546433333333
$ws.Range("N14").Value = 48.826393
$ws.Range("O14").Value = 0.06628560529319844
$ws.Range("P14").Value = 0.06628560529319844
$ws.Range("Q14").Value = 22.913005575075
$ws.Range("R14").Value = 206.217050175675
$ws.Range("S14").Value = 0.0002922451637972284
$ws.Range("T14").Value = 0.0002922451637972284
$ws.Range("G15").Value = 1.407825
$ws.Range("H15").Value = 4.223475
$ws.Range("I15").Value = 0.00440887825500804
$ws.Range("J15").Value = 0.00440887825500804
$ws.Range("O15").Value = 0.3480686258826592
$ws.Range("P15").Value = 0.3480686258826592
$ws.Range("Q15").Value = 120.317199037125
$ws.Range("R15").Value = 1082.854791334125
$ws.Range("S15").Value = 0.001534592195904585
$ws.Range("T15").Value = 0.001534592195904585
$ws.Range("G16").Value = 1.407825
$ws.Range("H16").Value = 4.223475
$ws.Range("I16").Value = 0.00440887825500804
$ws.Range("J16").Value = 0.00440887825500804
$ws.Range("M16").Value = 42.61351133333333
$ws.Range("N16").Value = 127.840534
$ws.Range("O16").Value = 0.17355341356458
$ws.Range("P16").Value = 0.17355341356458
$ws.Range("Q16").Value = 59.99236659284999
$ws.Range("R16").Value = 539.9312993356499
$ws.Range("S16").Value = 0.000765175871147294
$ws.Range("T16").Value = 0.0007651758711472943
$ws.Range("G17").Value = 1.407825
$ws.Range("H17").Value = 4.223475
$ws.Range("I17").Value = 0.00440887825500804
$ws.Range("J17").Value = 0.00440887825500804
$ws.Range("M17").Value = 101.183272
$ws.Range("N17").Value = 303.549816
$ws.Range("O17").Value = 0.4120923552595624
$ws.Range("P17").Value = 0.4120923552595624
$ws.Range("Q17").Value = 142.4483399034
$ws.Range("R17").Value = 1282.0350591306
$ws.Range("S17").Value = 0.001816865024158933
$ws.Range("T17").Value = 0.001816865024158933
